# Applies the perturbation-test edit described by the commit:
#   "Changed sheets in perturbation_tests/to_be_reformatted/math_L_curve"
#
# Net effect on the workbook:
#  - optimization_parameters sheet:
#      * the redundant extra header cells C1:F1 are removed
#      * row "Model" / "Sigmoid" becomes "production_function" / "Sigmoid"
#      * a new row "L_curve" / 1 is inserted right below it
#      * the old "Deletion" row (with values 0 / 3) is removed
#  - the optimization_parameters sheet tab becomes the active/selected one
#    (instead of optimization_diagnostics), with C1:F6 selected
#  - shared strings "Deletion" and "Model" drop out of use and are replaced
#    by new strings "production_function" and "L_curve" (handled
#    automatically by Excel when the workbook is saved)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# --- Row 1: drop the extra duplicated "value" header cells (C1:F1) ---
$ws.Range("C1:F1").ClearContents()

# --- Row 8: rename "Model" header to "production_function" ---
$ws.Cells.Item(8, 1).Value = "production_function"

# --- Insert a new row 9 for the "L_curve" parameter ---
$ws.Rows.Item(9).Insert()
$ws.Cells.Item(9, 1).Value = "L_curve"
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 2).NumberFormat = "0.00E+00"

# --- Remove the old "Deletion" row (now shifted down to row 17) ---
$ws.Rows.Item(17).Delete()

# --- Update view state: optimization_parameters becomes the active sheet ---
$ws.Select()
$ws.Range("C1:F6").Select()
